# Auto-generated script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.859.99"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "1.967.80"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "324.04"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "0.4783"
$ws.Range("E7").Value = "  -4.26%  "
$ws.Range("D8").Value = "0.4058"
$ws.Range("E8").Value = "  -4.35%  "
$ws.Range("D9").Value = "54.01"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "0.08547"
$ws.Range("E10").Value = "  -5.14%  "
$ws.Range("D11").Value = "1.064"
$ws.Range("E11").Value = "  -5.06%  "
$ws.Range("D12").Value = "22.54"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").Value = "1.957.93"
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("D14").Value = "7.687"
$ws.Range("E14").Value = "  -4.85%  "
$ws.Range("D15").Value = "6.266"
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("D16").Value = "1.012"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "90.23"
$ws.Range("E17").Value = "  -4.53%  "
$ws.Range("D18").Value = "'0.00001073"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("D19").Value = "0.06623"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "18.75"
$ws.Range("E20").Value = "  -5.33%  "
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "5.801"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").Value = "28.861.05"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").Value = "'11.60"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "2.196.95"
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("D27").Value = "154.41"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "20.28"
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("D29").Value = "5.983"
$ws.Range("E29").Value = "  -6.32%  "
$ws.Range("D30").Value = "2.161"
$ws.Range("E30").Value = "  -6.34%  "
$ws.Range("D31").Value = "124.54"
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("D32").Value = "1.011"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").Value = "0.09628"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").Value = "1.469"
$ws.Range("E34").Value = "  -6.40%  "
$ws.Range("D35").Value = "5.695"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").Value = "'3.690"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").Value = "0.02358"
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.273"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "8.812"
$ws.Range("E39").Value = "  -7.33%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06223"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").Value = "0.6269"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("D42").Value = "11.15"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").Value = "'1.010"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "0.1925"
$ws.Range("E44").Value = "  -6.49%  "
$ws.Range("D45").Value = "'1.350"
$ws.Range("E45").Value = "  +4.59%  "
$ws.Range("D46").Value = "0.5993"
$ws.Range("E46").Value = "  -5.58%  "
$ws.Range("D47").Value = "13.05"
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("D48").Value = "2.092"
$ws.Range("E48").Value = "  -5.22%  "
$ws.Range("D49").Value = "3.424"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06856"
$ws.Range("E51").Value = "  -1.88%  "
